$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update A283 (date) and D283 (low) per diff ---
$ws.Range("A283").Value() = 45449.2916666667
$ws.Range("D283").Value() = 6.09999990463257

# --- Add new row 284 (copy format of A283 so the date style matches) ---
$ws.Range("A283").Copy() | Out-Null
$ws.Range("A284").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A284").Value() = 45450.6060300926
$ws.Range("B284").Value() = 2600
$ws.Range("C284").Value() = 6.25
$ws.Range("D284").Value() = 6.09999990463257
$ws.Range("E284").Value() = 6.09999990463257
$ws.Range("F284").Value() = 6.25
$ws.Range("H284").Value() = $ws.Range("H283").Value()

# --- Rewrite the "adj_close" (G) column text values to match close (F) ---
# Use formulas producing text, then flatten via copy/paste-special so the
# cells remain plain shared-string text cells (matching original layout).
$ws.Range("G2").Formula = "=""8.53499984741211"""
$ws.Range("G3").Formula = "=""8.44999980926514"""
$ws.Range("G4").Formula = "=""8.47999954223633"""
$ws.Range("G5").Formula = "=""8.39999961853027"""
$ws.Range("G6").Formula = "=""8.42000007629395"""
$ws.Range("G7").Formula = "=""8.30000019073486"""
$ws.Range("G8").Formula = "=""8.23999977111816"""
$ws.Range("G9").Formula = "=""8.39999961853027"""
$ws.Range("G10").Formula = "=""8.39999961853027"""
$ws.Range("G11").Formula = "=""8.35999965667725"""
$ws.Range("G12").Formula = "=""8.31999969482422"""
$ws.Range("G13").Formula = "=""8.19999980926514"""
$ws.Range("G14").Formula = "=""8.22000026702881"""
$ws.Range("G15").Formula = "=""8.03999996185303"""
$ws.Range("G16").Formula = "=""8.11999988555908"""
$ws.Range("G17").Formula = "=""8.18000030517578"""
$ws.Range("G18").Formula = "=""8.15999984741211"""
$ws.Range("G19").Formula = "=""8.15999984741211"""
$ws.Range("G20").Formula = "=""8.19999980926514"""
$ws.Range("G21").Formula = "=""8.23999977111816"""
$ws.Range("G22").Formula = "=""8.23999977111816"""
$ws.Range("G23").Formula = "=""7.98000001907349"""
$ws.Range("G24").Formula = "=""7.80000019073486"""
$ws.Range("G25").Formula = "=""7.51999998092651"""
$ws.Range("G26").Formula = "=""8"""
$ws.Range("G27").Formula = "=""8.15999984741211"""
$ws.Range("G28").Formula = "=""8.15999984741211"""
$ws.Range("G29").Formula = "=""7.98000001907349"""
$ws.Range("G30").Formula = "=""7.98000001907349"""
$ws.Range("G31").Formula = "=""7.90000009536743"""
$ws.Range("G32").Formula = "=""7.90000009536743"""
$ws.Range("G33").Formula = "=""8.27999973297119"""
$ws.Range("G34").Formula = "=""8.38000011444092"""
$ws.Range("G35").Formula = "=""8.30000019073486"""
$ws.Range("G36").Formula = "=""8.38000011444092"""
$ws.Range("G37").Formula = "=""8.30000019073486"""
$ws.Range("G38").Formula = "=""8.39999961853027"""
$ws.Range("G39").Formula = "=""8.69999980926514"""
$ws.Range("G40").Formula = "=""8.52000045776367"""
$ws.Range("G41").Formula = "=""8.5"""
$ws.Range("G42").Formula = "=""8.39999961853027"""
$ws.Range("G43").Formula = "=""8.46000003814697"""
$ws.Range("G44").Formula = "=""8.68000030517578"""
$ws.Range("G45").Formula = "=""8.80000019073486"""
$ws.Range("G46").Formula = "=""8.84000015258789"""
$ws.Range("G47").Formula = "=""8.69999980926514"""
$ws.Range("G48").Formula = "=""9.0600004196167"""
$ws.Range("G49").Formula = "=""9.39999961853027"""
$ws.Range("G50").Formula = "=""9.72000026702881"""
$ws.Range("G51").Formula = "=""9.73999977111816"""
$ws.Range("G52").Formula = "=""9.69999980926514"""
$ws.Range("G53").Formula = "=""9.19999980926514"""
$ws.Range("G54").Formula = "=""9.02000045776367"""
$ws.Range("G55").Formula = "=""9.10000038146973"""
$ws.Range("G56").Formula = "=""9.22000026702881"""
$ws.Range("G57").Formula = "=""8.97999954223633"""
$ws.Range("G58").Formula = "=""8.89999961853027"""
$ws.Range("G59").Formula = "=""8.92000007629395"""
$ws.Range("G60").Formula = "=""8.80000019073486"""
$ws.Range("G61").Formula = "=""8.65999984741211"""
$ws.Range("G62").Formula = "=""8.5"""
$ws.Range("G63").Formula = "=""8.5"""
$ws.Range("G64").Formula = "=""8.5600004196167"""
$ws.Range("G65").Formula = "=""8.38000011444092"""
$ws.Range("G66").Formula = "=""8.30000019073486"""
$ws.Range("G67").Formula = "=""8.11999988555908"""
$ws.Range("G68").Formula = "=""8.19999980926514"""
$ws.Range("G69").Formula = "=""8.15999984741211"""
$ws.Range("G70").Formula = "=""8.07999992370605"""
$ws.Range("G71").Formula = "=""8"""
$ws.Range("G72").Formula = "=""8.15999984741211"""
$ws.Range("G73").Formula = "=""8.27999973297119"""
$ws.Range("G74").Formula = "=""8.27999973297119"""
$ws.Range("G75").Formula = "=""8.30000019073486"""
$ws.Range("G76").Formula = "=""7.96000003814697"""
$ws.Range("G77").Formula = "=""7.88000011444092"""
$ws.Range("G78").Formula = "=""7.92000007629395"""
$ws.Range("G79").Formula = "=""7.76000022888184"""
$ws.Range("G80").Formula = "=""7.90000009536743"""
$ws.Range("G81").Formula = "=""7.8600001335144"""
$ws.Range("G82").Formula = "=""7.8600001335144"""
$ws.Range("G83").Formula = "=""8"""
$ws.Range("G84").Formula = "=""8"""
$ws.Range("G85").Formula = "=""8.03999996185303"""
$ws.Range("G86").Formula = "=""8.03999996185303"""
$ws.Range("G87").Formula = "=""8.02000045776367"""
$ws.Range("G88").Formula = "=""8.03999996185303"""
$ws.Range("G89").Formula = "=""8"""
$ws.Range("G90").Formula = "=""8"""
$ws.Range("G91").Formula = "=""7.84000015258789"""
$ws.Range("G92").Formula = "=""8.10000038146973"""
$ws.Range("G93").Formula = "=""8.27999973297119"""
$ws.Range("G94").Formula = "=""8.10000038146973"""
$ws.Range("G95").Formula = "=""8.02000045776367"""
$ws.Range("G96").Formula = "=""8.02000045776367"""
$ws.Range("G97").Formula = "=""8.10000038146973"""
$ws.Range("G98").Formula = "=""8.4399995803833"""
$ws.Range("G99").Formula = "=""8.23999977111816"""
$ws.Range("G100").Formula = "=""8.23999977111816"""
$ws.Range("G101").Formula = "=""8.0600004196167"""
$ws.Range("G102").Formula = "=""8.19999980926514"""
$ws.Range("G103").Formula = "=""8.19999980926514"""
$ws.Range("G104").Formula = "=""8"""
$ws.Range("G105").Formula = "=""8.15999984741211"""
$ws.Range("G106").Formula = "=""8.07999992370605"""
$ws.Range("G107").Formula = "=""7.90000009536743"""
$ws.Range("G108").Formula = "=""7.90000009536743"""
$ws.Range("G109").Formula = "=""7.88000011444092"""
$ws.Range("G110").Formula = "=""7.67999982833862"""
$ws.Range("G111").Formula = "=""7.67999982833862"""
$ws.Range("G112").Formula = "=""7.82000017166138"""
$ws.Range("G113").Formula = "=""7.6399998664856"""
$ws.Range("G114").Formula = "=""7.6399998664856"""
$ws.Range("G115").Formula = "=""7.5"""
$ws.Range("G116").Formula = "=""7.53999996185303"""
$ws.Range("G117").Formula = "=""8.18000030517578"""
$ws.Range("G118").Formula = "=""8.23999977111816"""
$ws.Range("G119").Formula = "=""8.30000019073486"""
$ws.Range("G120").Formula = "=""8.30000019073486"""
$ws.Range("G121").Formula = "=""8.11999988555908"""
$ws.Range("G122").Formula = "=""8.11999988555908"""
$ws.Range("G123").Formula = "=""8.64000034332275"""
$ws.Range("G124").Formula = "=""8.26000022888184"""
$ws.Range("G125").Formula = "=""8.38000011444092"""
$ws.Range("G126").Formula = "=""8.19999980926514"""
$ws.Range("G127").Formula = "=""8.19999980926514"""
$ws.Range("G128").Formula = "=""8.19999980926514"""
$ws.Range("G129").Formula = "=""8.03999996185303"""
$ws.Range("G130").Formula = "=""8.03999996185303"""
$ws.Range("G131").Formula = "=""7.80000019073486"""
$ws.Range("G132").Formula = "=""7.80000019073486"""
$ws.Range("G133").Formula = "=""7.80000019073486"""
$ws.Range("G134").Formula = "=""7.96000003814697"""
$ws.Range("G135").Formula = "=""8.15999984741211"""
$ws.Range("G136").Formula = "=""8.15999984741211"""
$ws.Range("G137").Formula = "=""7.90000009536743"""
$ws.Range("G138").Formula = "=""7.90000009536743"""
$ws.Range("G139").Formula = "=""7.90000009536743"""
$ws.Range("G140").Formula = "=""7.88000011444092"""
$ws.Range("G141").Formula = "=""8.03999996185303"""
$ws.Range("G142").Formula = "=""8.03999996185303"""
$ws.Range("G143").Formula = "=""7.88000011444092"""
$ws.Range("G144").Formula = "=""7.69999980926514"""
$ws.Range("G145").Formula = "=""7.80000019073486"""
$ws.Range("G146").Formula = "=""7.80000019073486"""
$ws.Range("G147").Formula = "=""7.80000019073486"""
$ws.Range("G148").Formula = "=""7.94000005722046"""
$ws.Range("G149").Formula = "=""7.94000005722046"""
$ws.Range("G150").Formula = "=""7.6399998664856"""
$ws.Range("G151").Formula = "=""7.71999979019165"""
$ws.Range("G152").Formula = "=""7.55999994277954"""
$ws.Range("G153").Formula = "=""7.5"""
$ws.Range("G154").Formula = "=""7.65999984741211"""
$ws.Range("G155").Formula = "=""7.80000019073486"""
$ws.Range("G156").Formula = "=""7.80000019073486"""
$ws.Range("G157").Formula = "=""7.78000020980835"""
$ws.Range("G158").Formula = "=""7.5"""
$ws.Range("G159").Formula = "=""7.17999982833862"""
$ws.Range("G160").Formula = "=""7.09999990463257"""
$ws.Range("G161").Formula = "=""7.05999994277954"""
$ws.Range("G162").Formula = "=""7.05999994277954"""
$ws.Range("G163").Formula = "=""6.90000009536743"""
$ws.Range("G164").Formula = "=""7"""
$ws.Range("G165").Formula = "=""6.80000019073486"""
$ws.Range("G166").Formula = "=""7.15999984741211"""
$ws.Range("G167").Formula = "=""6.80000019073486"""
$ws.Range("G168").Formula = "=""6.59999990463257"""
$ws.Range("G169").Formula = "=""7"""
$ws.Range("G170").Formula = "=""7.19999980926514"""
$ws.Range("G171").Formula = "=""7.17999982833862"""
$ws.Range("G172").Formula = "=""7.15999984741211"""
$ws.Range("G173").Formula = "=""7.57999992370605"""
$ws.Range("G174").Formula = "=""7.30000019073486"""
$ws.Range("G175").Formula = "=""7.30000019073486"""
$ws.Range("G176").Formula = "=""7.30000019073486"""
$ws.Range("G177").Formula = "=""7.30000019073486"""
$ws.Range("G178").Formula = "=""7.28000020980835"""
$ws.Range("G179").Formula = "=""7.28000020980835"""
$ws.Range("G180").Formula = "=""7.01999998092651"""
$ws.Range("G181").Formula = "=""7"""
$ws.Range("G182").Formula = "=""6.92000007629395"""
$ws.Range("G183").Formula = "=""6.94000005722046"""
$ws.Range("G184").Formula = "=""7.19999980926514"""
$ws.Range("G185").Formula = "=""7.19999980926514"""
$ws.Range("G186").Formula = "=""7"""
$ws.Range("G187").Formula = "=""6.96000003814697"""
$ws.Range("G188").Formula = "=""7"""
$ws.Range("G189").Formula = "=""7"""
$ws.Range("G190").Formula = "=""6.98000001907349"""
$ws.Range("G191").Formula = "=""7.23999977111816"""
$ws.Range("G192").Formula = "=""7.40000009536743"""
$ws.Range("G193").Formula = "=""7.40000009536743"""
$ws.Range("G194").Formula = "=""7.40000009536743"""
$ws.Range("G195").Formula = "=""7.23999977111816"""
$ws.Range("G196").Formula = "=""7.44000005722046"""
$ws.Range("G197").Formula = "=""7.78000020980835"""
$ws.Range("G198").Formula = "=""7.59999990463257"""
$ws.Range("G199").Formula = "=""7.48000001907349"""
$ws.Range("G200").Formula = "=""7.32000017166138"""
$ws.Range("G201").Formula = "=""6.98000001907349"""
$ws.Range("G202").Formula = "=""6.71999979019165"""
$ws.Range("G203").Formula = "=""6.71999979019165"""
$ws.Range("G204").Formula = "=""6.71999979019165"""
$ws.Range("G205").Formula = "=""6.57999992370605"""
$ws.Range("G206").Formula = "=""6.44000005722046"""
$ws.Range("G207").Formula = "=""6.46000003814697"""
$ws.Range("G208").Formula = "=""6.48000001907349"""
$ws.Range("G209").Formula = "=""6.90000009536743"""
$ws.Range("G210").Formula = "=""6.59999990463257"""
$ws.Range("G211").Formula = "=""6.5"""
$ws.Range("G212").Formula = "=""6.40000009536743"""
$ws.Range("G213").Formula = "=""6.26000022888184"""
$ws.Range("G214").Formula = "=""6.1399998664856"""
$ws.Range("G215").Formula = "=""6.07999992370605"""
$ws.Range("G216").Formula = "=""6.07999992370605"""
$ws.Range("G217").Formula = "=""6.17999982833862"""
$ws.Range("G218").Formula = "=""6.17999982833862"""
$ws.Range("G219").Formula = "=""6.09999990463257"""
$ws.Range("G220").Formula = "=""6.09999990463257"""
$ws.Range("G221").Formula = "=""6.32000017166138"""
$ws.Range("G222").Formula = "=""6.17999982833862"""
$ws.Range("G223").Formula = "=""6.17999982833862"""
$ws.Range("G224").Formula = "=""6.17999982833862"""
$ws.Range("G225").Formula = "=""6.09999990463257"""
$ws.Range("G226").Formula = "=""6"""
$ws.Range("G227").Formula = "=""5.90000009536743"""
$ws.Range("G228").Formula = "=""5.88000011444092"""
$ws.Range("G229").Formula = "=""5.78000020980835"""
$ws.Range("G230").Formula = "=""5.80000019073486"""
$ws.Range("G231").Formula = "=""5.80000019073486"""
$ws.Range("G232").Formula = "=""5.90000009536743"""
$ws.Range("G233").Formula = "=""5.90000009536743"""
$ws.Range("G234").Formula = "=""6"""
$ws.Range("G235").Formula = "=""5.90000009536743"""
$ws.Range("G236").Formula = "=""5.90000009536743"""
$ws.Range("G237").Formula = "=""5.75"""
$ws.Range("G238").Formula = "=""5.65000009536743"""
$ws.Range("G239").Formula = "=""5.80000019073486"""
$ws.Range("G240").Formula = "=""5.84999990463257"""
$ws.Range("G241").Formula = "=""5.80000019073486"""
$ws.Range("G242").Formula = "=""5.80000019073486"""
$ws.Range("G243").Formula = "=""5.75"""
$ws.Range("G244").Formula = "=""6.30000019073486"""
$ws.Range("G245").Formula = "=""6.19999980926514"""
$ws.Range("G246").Formula = "=""6.05000019073486"""
$ws.Range("G247").Formula = "=""6.05000019073486"""
$ws.Range("G248").Formula = "=""6.05000019073486"""
$ws.Range("G249").Formula = "=""5.90000009536743"""
$ws.Range("G250").Formula = "=""5.75"""
$ws.Range("G251").Formula = "=""5.5"""
$ws.Range("G252").Formula = "=""5.34999990463257"""
$ws.Range("G253").Formula = "=""5.19999980926514"""
$ws.Range("G254").Formula = "=""5.05000019073486"""
$ws.Range("G255").Formula = "=""5.34999990463257"""
$ws.Range("G256").Formula = "=""5.34999990463257"""
$ws.Range("G257").Formula = "=""5.44999980926514"""
$ws.Range("G258").Formula = "=""5.55000019073486"""
$ws.Range("G259").Formula = "=""5.5"""
$ws.Range("G260").Formula = "=""5.25"""
$ws.Range("G261").Formula = "=""5"""
$ws.Range("G262").Formula = "=""5"""
$ws.Range("G263").Formula = "=""5.40000009536743"""
$ws.Range("G264").Formula = "=""5.30000019073486"""
$ws.Range("G265").Formula = "=""5.25"""
$ws.Range("G266").Formula = "=""5.34999990463257"""
$ws.Range("G267").Formula = "=""5.40000009536743"""
$ws.Range("G268").Formula = "=""5.59999990463257"""
$ws.Range("G269").Formula = "=""5.80000019073486"""
$ws.Range("G270").Formula = "=""5.65000009536743"""
$ws.Range("G271").Formula = "=""5.59999990463257"""
$ws.Range("G272").Formula = "=""5.59999990463257"""
$ws.Range("G273").Formula = "=""5.80000019073486"""
$ws.Range("G274").Formula = "=""5.69999980926514"""
$ws.Range("G275").Formula = "=""5.80000019073486"""
$ws.Range("G276").Formula = "=""6.05000019073486"""
$ws.Range("G277").Formula = "=""6.05000019073486"""
$ws.Range("G278").Formula = "=""6"""
$ws.Range("G279").Formula = "=""6"""
$ws.Range("G280").Formula = "=""6"""
$ws.Range("G281").Formula = "=""5.80000019073486"""
$ws.Range("G282").Formula = "=""6.05000019073486"""
$ws.Range("G283").Formula = "=""6.09999990463257"""
$ws.Range("G284").Formula = "=""6.25"""

$gRange = $ws.Range("G2:G284")
$gRange.Copy() | Out-Null
$gRange.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

Write-Host "done"